$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the numeric results for Mississippi (row 42) back to blank text cells,
# mirroring a failed scrape run: columns B..H become empty text values and
# lose any previously-applied formatting (e.g. B42's date format).
$blankCells = @("B42", "C42", "D42", "E42", "F42", "G42", "H42")
foreach ($addr in $blankCells) {
    $rng = $ws.Range($addr)
    $rng.Value = "'"
    $rng.ClearFormats()
}

# The disparity flags also reset to False for this failed run.
$ws.Range("I42").Value = $false
$ws.Range("J42").Value = $false

# Record the connection error that caused the run to fail.
$ws.Range("O42").Value = "An error occurred. ... ConnectionError(ProtocolError('Connection aborted.', ConnectionResetError(104, 'Connection reset by peer')))"
